$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unprotect the sheet so the data cells (currently locked) can be edited.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure text (A10).
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-04 for illustrative purposes only and are subject to change."
# Re-fit the row height (editing the multi-line text can otherwise stamp an
# explicit row height override that wasn't in the original file).
$ws.Rows(10).AutoFit()

# Refresh the Weight / Percent Change figures.
$ws.Range("D2").Value = 0.249488402130182
$ws.Range("E2").Value = -0.01394585726004915

$ws.Range("D3").Value = 0.4934898340780004
$ws.Range("E3").Value = 0.002917771883289122

$ws.Range("D4").Value = 0.09892064364316007
$ws.Range("E4").Value = -0.0170400153168675

$ws.Range("D5").Value = 0.1008559221016501
$ws.Range("E5").Value = 0.003555798687089551

$ws.Range("D6").Value = 0.05724519804700743
$ws.Range("E6").Value = -0.009143621766280252

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = -0.003889853247900588

# Restore worksheet protection to match the original state.
$ws.Protect("wb'm|;4", $true, $true, $true, $false)
